# Update DAMSLTag (col I) and DialogAct (col J) values per the dialog-act re-annotation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=4; I='sd'; J='Statement-non-opinion'},
    @{Row=8; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=10; I='aa'; J='Agree/Accept'},
    @{Row=12; I='sd'; J='Statement-non-opinion'},
    @{Row=18; I='sv'; J='Statement-opinion'},
    @{Row=20; I='sd'; J='Statement-non-opinion'},
    @{Row=21; I='sv'; J='Statement-opinion'},
    @{Row=25; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=29; I='ba'; J='Appreciation'},
    @{Row=46; I='%'; J='Uninterpretable'},
    @{Row=48; I='sd'; J='Statement-non-opinion'},
    @{Row=66; I='aa'; J='Agree/Accept'},
    @{Row=67; I='aa'; J='Agree/Accept'},
    @{Row=70; I='sd'; J='Statement-non-opinion'},
    @{Row=84; I='sd'; J='Statement-non-opinion'},
    @{Row=92; I='sd'; J='Statement-non-opinion'},
    @{Row=108; I='sd'; J='Statement-non-opinion'},
    @{Row=112; I='aa'; J='Agree/Accept'},
    @{Row=114; I='sd'; J='Statement-non-opinion'},
    @{Row=124; I='aa'; J='Agree/Accept'},
    @{Row=127; I='sv'; J='Statement-opinion'},
    @{Row=142; I='sd'; J='Statement-non-opinion'},
    @{Row=143; I='sd'; J='Statement-non-opinion'},
    @{Row=145; I='sd'; J='Statement-non-opinion'},
    @{Row=148; I='ba'; J='Appreciation'},
    @{Row=154; I='sv'; J='Statement-opinion'},
    @{Row=169; I='sd'; J='Statement-non-opinion'},
    @{Row=173; I='ba'; J='Appreciation'},
    @{Row=176; I='ba'; J='Appreciation'},
    @{Row=193; I='sv'; J='Statement-opinion'},
    @{Row=194; I='ba'; J='Appreciation'},
    @{Row=199; I='qy'; J='Yes-No-Question'},
    @{Row=200; I='aa'; J='Agree/Accept'},
    @{Row=201; I='sv'; J='Statement-opinion'},
    @{Row=204; I='aa'; J='Agree/Accept'},
    @{Row=212; I='sv'; J='Statement-opinion'},
    @{Row=213; I='sd'; J='Statement-non-opinion'},
    @{Row=228; I='sv'; J='Statement-opinion'},
    @{Row=230; I='sv'; J='Statement-opinion'},
    @{Row=239; I='aa'; J='Agree/Accept'},
    @{Row=255; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=256; I='sv'; J='Statement-opinion'},
    @{Row=259; I='sd'; J='Statement-non-opinion'},
    @{Row=275; I='aa'; J='Agree/Accept'},
    @{Row=287; I='sd'; J='Statement-non-opinion'},
    @{Row=293; I='aa'; J='Agree/Accept'},
    @{Row=298; I='sd'; J='Statement-non-opinion'},
    @{Row=300; I='sv'; J='Statement-opinion'},
    @{Row=305; I='ba'; J='Appreciation'},
    @{Row=312; I='sd'; J='Statement-non-opinion'},
    @{Row=321; I='sv'; J='Statement-opinion'},
    @{Row=339; I='qy'; J='Yes-No-Question'},
    @{Row=340; I='sd'; J='Statement-non-opinion'},
    @{Row=352; I='aa'; J='Agree/Accept'},
    @{Row=359; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=361; I='sd'; J='Statement-non-opinion'},
    @{Row=378; I='sd'; J='Statement-non-opinion'},
    @{Row=379; I='aa'; J='Agree/Accept'},
    @{Row=383; I='sd'; J='Statement-non-opinion'},
    @{Row=396; I='sd'; J='Statement-non-opinion'},
    @{Row=398; I='sd'; J='Statement-non-opinion'},
    @{Row=409; I='sd'; J='Statement-non-opinion'},
    @{Row=414; I='aa'; J='Agree/Accept'},
    @{Row=423; I='qy'; J='Yes-No-Question'},
    @{Row=425; I='sd'; J='Statement-non-opinion'},
    @{Row=431; I='sd'; J='Statement-non-opinion'},
    @{Row=441; I='sd'; J='Statement-non-opinion'},
    @{Row=446; I='sd'; J='Statement-non-opinion'},
    @{Row=447; I='aa'; J='Agree/Accept'},
    @{Row=450; I='sd'; J='Statement-non-opinion'},
    @{Row=451; I='sd'; J='Statement-non-opinion'},
    @{Row=454; I='sv'; J='Statement-opinion'},
    @{Row=455; I='sd'; J='Statement-non-opinion'},
    @{Row=457; I='sd'; J='Statement-non-opinion'},
    @{Row=461; I='ba'; J='Appreciation'},
    @{Row=474; I='sd'; J='Statement-non-opinion'},
    @{Row=479; I='sd'; J='Statement-non-opinion'},
    @{Row=482; I='sd'; J='Statement-non-opinion'},
    @{Row=483; I='sd'; J='Statement-non-opinion'},
    @{Row=489; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=490; I='sd'; J='Statement-non-opinion'},
    @{Row=492; I='sv'; J='Statement-opinion'},
    @{Row=498; I='sv'; J='Statement-opinion'},
    @{Row=502; I='sd'; J='Statement-non-opinion'},
    @{Row=504; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows."
